# Apply the "exp1. with proper decay group!" edit:
# populate the previously-empty "charge" sheet with a header row and one
# data row describing a new experiment run.

$wb = $excel.ActiveWorkbook
$expWs = $wb.Worksheets.Item("exp")
$ws = $wb.Worksheets.Item("charge")

# --- Write cells in the exact order that reproduces the shared-string table ---
# (new strings must be introduced in this order: git hash, run name, then the
#  four new header labels, then the optimizer value)
$ws.Range("V2").Value = "c62429ece4a9ac7417a635f5932dae136bdd8f6d"
$ws.Range("A2").Value = "v20240628210409"
$ws.Range("N1").Value = "optimizer"
$ws.Range("F1").Value = "epochs"
$ws.Range("O1").Value = "reinitialize_n_layers_of_backbone"
$ws.Range("P1").Value = "llrd"
$ws.Range("N2").Value = "adamw-forgot-wd"

# --- Header row (reusing strings already present from the "exp" sheet) ---
$ws.Range("A1").Value = "run name"
$ws.Range("B1").Value = "external"
$ws.Range("C1").Value = "resampling"
$ws.Range("D1").Value = "seed"
$ws.Range("E1").Value = "effective bs"
$ws.Range("G1").Value = "effective lr"
$ws.Range("H1").Value = "scheduler"
$ws.Range("I1").Value = "warmup"
$ws.Range("J1").Value = "model"
$ws.Range("K1").Value = "task"
$ws.Range("L1").Value = "pool"
$ws.Range("M1").Value = "loss"
$ws.Range("Q1").Value = "threshold"
$ws.Range("R1").Value = "train/infer max len"
$ws.Range("S1").Value = "oof"
$ws.Range("T1").Value = "lb"
$ws.Range("U1").Value = "diff"
$ws.Range("V1").Value = "git"
$ws.Range("W1").Value = "kaggle"

# --- Data row 2 ---
$ws.Range("B2").Value = "yes"
$ws.Range("C2").Value = "skf-5-42"
$ws.Range("D2").Value = 20230310
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = "'1.00E-05"
$ws.Range("H2").Value = "linear"
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = "deberta small"
$ws.Range("K2").Value = "reg"
$ws.Range("L2").Value = "null"
$ws.Range("M2").Value = "mse"
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = "null"
$ws.Range("Q2").Value = "no"
$ws.Range("R2").Value = "1024/1024"
$ws.Range("W2").Value = "https://www.kaggle.com/code/reighns/v20240624231703-all-folds-inference?scriptVersionId=185266778"

# --- Column widths (best-fit similar to the "exp" sheet's layout) ---
$ws.Columns.Item(1).ColumnWidth = 16
$ws.Columns.Item(2).ColumnWidth = 7.83203125
$ws.Columns.Item(3).ColumnWidth = 10
$ws.Columns.Item(4).ColumnWidth = 9.1640625
$ws.Columns.Item(5).ColumnWidth = 10.33203125
$ws.Columns.Item(6).ColumnWidth = 10.33203125
$ws.Columns.Item(7).ColumnWidth = 9.6640625
$ws.Columns.Item(8).ColumnWidth = 9.1640625
$ws.Columns.Item(9).ColumnWidth = 7.6640625
$ws.Columns.Item(10).ColumnWidth = 12.1640625
$ws.Columns.Item(11).ColumnWidth = 4.5
$ws.Columns.Item(12).ColumnWidth = 4.6640625
$ws.Columns.Item(13).ColumnWidth = 4.5
$ws.Columns.Item(14).ColumnWidth = 8.83203125
$ws.Columns.Item(15).ColumnWidth = 28.83203125
$ws.Columns.Item(16).ColumnWidth = 4.1640625
$ws.Columns.Item(17).ColumnWidth = 8.83203125
$ws.Columns.Item(18).ColumnWidth = 15.83203125
$ws.Columns.Item(19).ColumnWidth = 12.1640625
$ws.Columns.Item(20).ColumnWidth = 6.1640625
$ws.Columns.Item(21).ColumnWidth = 12.1640625
$ws.Columns.Item(22).ColumnWidth = 41
$ws.Columns.Item(23).ColumnWidth = 88.1640625

# --- Selection state ---
$ws.Range("N3").Select() | Out-Null
$expWs.Range("A2:S2").Select() | Out-Null
$ws.Activate() | Out-Null
